$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for all existing data rows (2-119)
for ($r = 2; $r -le 119; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 119 gains an explicit row height (15, custom)
$ws.Rows.Item(119).RowHeight = 15

# Copy formatting from row 119 into the new row 120, then fill in its values
$ws.Range("A119:R119").Copy() | Out-Null
$ws.Range("A120:R120").PasteSpecial(-4122) | Out-Null
$ws.Range("F120").ClearContents() | Out-Null

$ws.Range("A120").Value = "A 46757-2023"
$ws.Range("B120").Value = 45198
$ws.Range("C120").Value = 45202
$ws.Range("D120").Value = "STOCKHOLMS LÄN"
$ws.Range("E120").Value = "SIGTUNA"
$ws.Range("G120").Value = 5.3
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = 0
$ws.Range("N120").Value = 0
$ws.Range("O120").Value = 0
$ws.Range("P120").Value = 0
$ws.Range("Q120").Value = 0
$ws.Range("R120").Value = ""

Write-Host "Update complete"
